$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.979788333333333
$ws.Range("H2").Value = 14.939365
$ws.Range("I2").Value = 0.129176854764059
$ws.Range("J2").Value = 0.129176854764059
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 203.0691603333333
$ws.Range("N2").Value = 609.207481
$ws.Range("O2").Value = 0.9796789863919257
$ws.Range("P2").Value = 0.9796789863919257
$ws.Range("Q2").Value = 1011.24143548773
$ws.Range("R2").Value = 9101.172919389566
$ws.Range("S2").Value = 0.1265518501405503
$ws.Range("T2").Value = 0.1265518501405503
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.979788333333333
$ws.Range("H3").Value = 14.939365
$ws.Range("I3").Value = 0.129176854764059
$ws.Range("J3").Value = 0.129176854764059
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.5294913333333333
$ws.Range("N3").Value = 1.588474
$ws.Range("O3").Value = 0.002554457466076205
$ws.Range("P3").Value = 0.002554457466076205
$ws.Range("Q3").Value = 2.636754764334444
$ws.Range("R3").Value = 23.73079287901
$ws.Range("S3").Value = 0.000329976781096292
$ws.Range("T3").Value = 0.0003299767810962921
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.979788333333333
$ws.Range("H4").Value = 14.939365
$ws.Range("I4").Value = 0.129176854764059
$ws.Range("J4").Value = 0.129176854764059
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.1824963333333333
$ws.Range("N4").Value = 0.547489
$ws.Range("O4").Value = 0.0008804282371915408
$ws.Range("P4").Value = 0.0008804282371915407
$ws.Range("Q4").Value = 0.9087931116094445
$ws.Range("R4").Value = 8.179138004485001
$ws.Range("S4").Value = 0.0001137309505258681
$ws.Range("T4").Value = 0.0001137309505258681
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 4.979788333333333
$ws.Range("H5").Value = 14.939365
$ws.Range("I5").Value = 0.129176854764059
$ws.Range("J5").Value = 0.129176854764059
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.500179
$ws.Range("N5").Value = 10.500537
$ws.Range("O5").Value = 0.01688612790480639
$ws.Range("P5").Value = 0.01688612790480639
$ws.Range("Q5").Value = 17.43015054877833
$ws.Range("R5").Value = 156.871354939005
$ws.Range("S5").Value = 0.002181296891886499
$ws.Range("T5").Value = 0.002181296891886499
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 18.019504
$ws.Range("H6").Value = 54.058512
$ws.Range("I6").Value = 0.467430078412646
$ws.Range("J6").Value = 0.4674300784126461
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 203.0691603333333
$ws.Range("N6").Value = 609.207481
$ws.Range("O6").Value = 0.9796789863919257
$ws.Range("P6").Value = 0.9796789863919257
$ws.Range("Q6").Value = 3659.205546903142
$ws.Range("R6").Value = 32932.84992212828
$ws.Range("S6").Value = 0.4579314254283994
$ws.Range("T6").Value = 0.4579314254283994
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 18.019504
$ws.Range("H7").Value = 54.058512
$ws.Range("I7").Value = 0.467430078412646
$ws.Range("J7").Value = 0.4674300784126461
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.5294913333333333
$ws.Range("N7").Value = 1.588474
$ws.Range("O7").Value = 0.002554457466076205
$ws.Range("P7").Value = 0.002554457466076205
$ws.Range("Q7").Value = 9.541171198965333
$ws.Range("R7").Value = 85.870540790688
$ws.Range("S7").Value = 0.001194030253669769
$ws.Range("T7").Value = 0.00119403025366977
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 18.019504
$ws.Range("H8").Value = 54.058512
$ws.Range("I8").Value = 0.467430078412646
$ws.Range("J8").Value = 0.4674300784126461
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.1824963333333333
$ws.Range("N8").Value = 0.547489
$ws.Range("O8").Value = 0.0008804282371915408
$ws.Range("P8").Value = 0.0008804282371915407
$ws.Range("Q8").Value = 3.288493408485334
$ws.Range("R8").Value = 29.596440676368
$ws.Range("S8").Value = 0.0004115386399471496
$ws.Range("T8").Value = 0.0004115386399471496
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 18.019504
$ws.Range("H9").Value = 54.058512
$ws.Range("I9").Value = 0.467430078412646
$ws.Range("J9").Value = 0.4674300784126461
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.500179
$ws.Range("N9").Value = 10.500537
$ws.Range("O9").Value = 0.01688612790480639
$ws.Range("P9").Value = 0.01688612790480639
$ws.Range("Q9").Value = 63.071489491216
$ws.Range("R9").Value = 567.643405420944
$ws.Range("S9").Value = 0.007893084090629622
$ws.Range("T9").Value = 0.007893084090629623
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 8.752692000000001
$ws.Range("H10").Value = 26.258076
$ws.Range("I10").Value = 0.2270468436801446
$ws.Range("J10").Value = 0.2270468436801446
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 203.0691603333333
$ws.Range("N10").Value = 609.207481
$ws.Range("O10").Value = 0.9796789863919257
$ws.Range("P10").Value = 0.9796789863919257
$ws.Range("Q10").Value = 1777.401815096284
$ws.Range("R10").Value = 15996.61633586656
$ws.Range("S10").Value = 0.2224330216800501
$ws.Range("T10").Value = 0.2224330216800501
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 8.752692000000001
$ws.Range("H11").Value = 26.258076
$ws.Range("I11").Value = 0.2270468436801446
$ws.Range("J11").Value = 0.2270468436801446
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.5294913333333333
$ws.Range("N11").Value = 1.588474
$ws.Range("O11").Value = 0.002554457466076205
$ws.Range("P11").Value = 0.002554457466076205
$ws.Range("Q11").Value = 4.634474557336
$ws.Range("R11").Value = 41.710271016024
$ws.Range("S11").Value = 0.0005799815049877823
$ws.Range("T11").Value = 0.0005799815049877823
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 8.752692000000001
$ws.Range("H12").Value = 26.258076
$ws.Range("I12").Value = 0.2270468436801446
$ws.Range("J12").Value = 0.2270468436801446
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.1824963333333333
$ws.Range("N12").Value = 0.547489
$ws.Range("O12").Value = 0.0008804282371915408
$ws.Range("P12").Value = 0.0008804282371915407
$ws.Range("Q12").Value = 1.597334196796
$ws.Range("R12").Value = 14.376007771164
$ws.Range("S12").Value = 0.000199898452341213
$ws.Range("T12").Value = 0.000199898452341213
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 8.752692000000001
$ws.Range("H13").Value = 26.258076
$ws.Range("I13").Value = 0.2270468436801446
$ws.Range("J13").Value = 0.2270468436801446
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 3.500179
$ws.Range("N13").Value = 10.500537
$ws.Range("O13").Value = 0.01688612790480639
$ws.Range("P13").Value = 0.01688612790480639
$ws.Range("Q13").Value = 30.635988731868
$ws.Range("R13").Value = 275.723898586812
$ws.Range("S13").Value = 0.003833942042765505
$ws.Range("T13").Value = 0.003833942042765505
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 6.798175000000001
$ws.Range("H14").Value = 20.394525
$ws.Range("I14").Value = 0.1763462231431503
$ws.Range("J14").Value = 0.1763462231431503
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 203.0691603333333
$ws.Range("N14").Value = 609.207481
$ws.Range("O14").Value = 0.9796789863919257
$ws.Range("P14").Value = 0.9796789863919257
$ws.Range("Q14").Value = 1380.499689049059
$ws.Range("R14").Value = 12424.49720144153
$ws.Range("S14").Value = 0.1727626891429259
$ws.Range("T14").Value = 0.1727626891429259
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 6.798175000000001
$ws.Range("H15").Value = 20.394525
$ws.Range("I15").Value = 0.1763462231431503
$ws.Range("J15").Value = 0.1763462231431503
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 0.5294913333333333
$ws.Range("N15").Value = 1.588474
$ws.Range("O15").Value = 0.002554457466076205
$ws.Range("P15").Value = 0.002554457466076205
$ws.Range("Q15").Value = 3.599574744983333
$ws.Range("R15").Value = 32.39617270485
$ws.Range("S15").Value = 0.0004504689263223608
$ws.Range("T15").Value = 0.0004504689263223609
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 6.798175000000001
$ws.Range("H16").Value = 20.394525
$ws.Range("I16").Value = 0.1763462231431503
$ws.Range("J16").Value = 0.1763462231431503
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.1824963333333333
$ws.Range("N16").Value = 0.547489
$ws.Range("O16").Value = 0.0008804282371915408
$ws.Range("P16").Value = 0.0008804282371915407
$ws.Range("Q16").Value = 1.240642010858334
$ws.Range("R16").Value = 11.165778097725
$ws.Range("S16").Value = 0.0001552601943773099
$ws.Range("T16").Value = 0.0001552601943773099
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 6.798175000000001
$ws.Range("H17").Value = 20.394525
$ws.Range("I17").Value = 0.1763462231431503
$ws.Range("J17").Value = 0.1763462231431503
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 3.500179
$ws.Range("N17").Value = 10.500537
$ws.Range("O17").Value = 0.01688612790480639
$ws.Range("P17").Value = 0.01688612790480639
$ws.Range("Q17").Value = 17.43015054877833
$ws.Range("R17").Value = 214.153464359925
$ws.Range("S17").Value = 0.002977804879524766
$ws.Range("T17").Value = 0.002977804879524766
